$d = $word.ActiveDocument

function Find-ParagraphIndex($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$substr*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "March " + "1" + "9" + ", 2024"  ->  "March " + "19" + ", 2024"
#    (the two single-character runs "1" and "9" become one "19" run)
# ---------------------------------------------------------------------------
$findRng = $d.Content
$foundHeading = $findRng.Find.Execute("March 19, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundHeading) {
    throw "Could not find the 'March 19, 2024' heading paragraph"
}

$mergeRng = $d.Content
$mergedOk = $mergeRng.Find.Execute("19", $true, $false, $false, $false, $false, $true, 1, $false, "19", 2)
if (-not $mergedOk) {
    throw "Could not merge the '1' and '9' runs into '19'"
}

# ---------------------------------------------------------------------------
# 2) Insert the new March 20, 2024 daily-log block right after the paragraph
#    that ends with "...NN_implementation jupyter Notebook used in VM Ubuntu."
#    and right before the existing blank paragraph that follows it.
# ---------------------------------------------------------------------------
$anchorIdx = Find-ParagraphIndex("NN_implementation jupyter Notebook used in VM Ubuntu.")
if ($anchorIdx -eq -1) {
    throw "Could not find anchor paragraph (NN_implementation...)"
}

# Keep formatting templates (their own paragraph mark included) before any
# insertions shift paragraph indices around.
$tnrTemplateRange = $d.Paragraphs.Item($anchorIdx).Range

# --- 2a) blank (Times New Roman) paragraph right after the anchor ----------
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$insPoint1 = $anchorPara.Range.Duplicate
$insPoint1.Collapse(0)
$insPoint1.InsertParagraphAfter()
$blankIdx = $anchorIdx + 1

# --- 2b) "March 20, 2024" heading paragraph (Calibri, Bold, 24) ------------
#         built by cloning the existing "March 19, 2024" heading (so it
#         picks up the exact same run formatting), then patching the date.
$headingIdx = Find-ParagraphIndex("March 19, 2024")
if ($headingIdx -eq -1) {
    throw "Could not re-locate the 'March 19, 2024' heading paragraph"
}
$headingTemplateRange = $d.Paragraphs.Item($headingIdx).Range

$newHeadingInsertStart = $d.Paragraphs.Item($blankIdx + 1).Range.Start
$newHeadingDest = $d.Range($newHeadingInsertStart, $newHeadingInsertStart)
$newHeadingDest.FormattedText = $headingTemplateRange.FormattedText
$newHeadingIdx = $blankIdx + 1

$newHeadingStart = $d.Paragraphs.Item($newHeadingIdx).Range.Start
$dateDigitsRng = $d.Range($newHeadingStart + 6, $newHeadingStart + 8)
if ($dateDigitsRng.Text -ne "19") {
    throw "Unexpected text while patching the new heading date: [$($dateDigitsRng.Text)]"
}
$dateDigitsRng.Text = "20"

# --- 2c) daily-log content paragraph (Times New Roman, 20) -----------------
#         built by cloning the existing body paragraph's formatting, then
#         replacing its text.
$newHeadingPara = $d.Paragraphs.Item($newHeadingIdx)
$contentInsertPos = $newHeadingPara.Range.End
$contentDest = $d.Range($contentInsertPos, $contentInsertPos)
$contentDest.FormattedText = $tnrTemplateRange.FormattedText
$contentIdx = $newHeadingIdx + 1

$contentPara = $d.Paragraphs.Item($contentIdx)
$contentStart = $contentPara.Range.Start
$contentEnd = $contentPara.Range.End - 1
$contentTextRng = $d.Range($contentStart, $contentEnd)

$contentText = "Implementation Hadoop part completed, also dataset" + `
    " origin is mentioned and how size was increased. " + `
    "Uploading " + `
    "1.Increasing_dataset_size.ipynb" + `
    "." + `
    " " + `
    "Starting with MySQL dataload."
$contentTextRng.Text = $contentText

# --- 2d) trailing blank (Times New Roman) paragraph ------------------------
$contentParaAgain = $d.Paragraphs.Item($contentIdx)
$insPoint2 = $contentParaAgain.Range.Duplicate
$insPoint2.Collapse(0)
$insPoint2.InsertParagraphAfter()

"Edit applied successfully"
